# Add "Romania" and "Slovakia" market test-data sheets, cloned from the
# existing "Belgium" sheet (same layout/styles), and update the UK sheet's
# selection state + which tab is active - matching the authored diff.

$wb = $excel.ActiveWorkbook

# --- UK sheet: select whole column A, no longer the active tab -------------
$uk = $wb.Worksheets.Item("UK")
$uk.Activate()
$uk.Range("A1:XFD1048576").Select() | Out-Null

# --- Romania: clone of Belgium with Romania-specific market/story values ---
$belgium = $wb.Worksheets.Item("Belgium")

$belgium.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count)) | Out-Null
$romania = $wb.Worksheets.Item($wb.Worksheets.Count)
$romania.Name = "Romania"
$romania.Range("B2").Value = "Romania Market"
$romania.Range("B4").Value = "NGC-4307/T3536/T3543"
$romania.Range("B4").Select() | Out-Null

# --- Slovakia: clone of Belgium, left blank (market/story not filled yet) --
$belgium.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count)) | Out-Null
$slovakia = $wb.Worksheets.Item($wb.Worksheets.Count)
$slovakia.Name = "Slovakia"
$slovakia.Range("B2").ClearContents() | Out-Null
$slovakia.Range("B4").ClearContents() | Out-Null
$slovakia.Range("A9").Select() | Out-Null

# Romania ends up as the active tab.
$romania.Activate() | Out-Null
